$d = $word.ActiveDocument

# "Se levantan Faltas Leves: "  ->  "Se levantan Defectos leves: "
$d.Content.Find.Execute("Se levantan Faltas Leves: ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Se levantan Defectos leves: ", 2)

# "Se Mantienen Faltas Leves: "  ->  "Se Mantienen Defectos Leves: "
$d.Content.Find.Execute("Se Mantienen Faltas Leves: ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Se Mantienen Defectos Leves: ", 2)

# "FALTAS" (table header) -> "Defectos"
$d.Content.Find.Execute("FALTAS", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Defectos", 2)
